# Auto-generated from the commit diff: updates D (Price) / E (Volume(1h))
# text cells in cryptos.xlsx, plus the Toncoin/Stellar row swap (rows 26-27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.778.22"
$ws.Range("E2").Value = "  +6.46%  "
$ws.Range("D3").Value = "1.737.33"
$ws.Range("E3").Value = "  +5.19%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'227.37"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").Value = "'0.5458"
$ws.Range("E6").Value = "  +4.03%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.2762"
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").Value = "'0.06718"
$ws.Range("E9").Value = "  +5.53%  "
$ws.Range("D10").Value = "'21.96"
$ws.Range("E10").Value = "  +7.04%  "
$ws.Range("D11").Value = "'0.07775"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "'4.683"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.741.02"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "1.977.31"
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "'0.5968"
$ws.Range("E15").Value = "  +6.23%  "
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "'69.05"
$ws.Range("E17").Value = "  +5.52%  "
$ws.Range("D18").Value = "27.776.75"
$ws.Range("E18").Value = "  +6.46%  "
$ws.Range("D19").Value = "'225.50"
$ws.Range("E19").Value = "  +18.13%  "
$ws.Range("D20").Value = "'4.823"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'10.89"
$ws.Range("E22").Value = "  +5.40%  "
$ws.Range("D23").Value = "'6.220"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "'147.19"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1249"
$ws.Range("E26").Value = "  +4.04%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.702"
$ws.Range("E27").Value = "  +13.66%  "
$ws.Range("D28").Value = "'7.451"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").Value = "'17.12"
$ws.Range("E29").Value = "  +7.38%  "
$ws.Range("D30").Value = "'0.05662"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "'1.311"
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("D32").Value = "'3.691"
$ws.Range("E32").Value = "  +5.76%  "
$ws.Range("D33").Value = "'3.509"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").Value = "'1.678"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("D35").Value = "'0.9767"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").Value = "'2.853"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").Value = "'2.451"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").Value = "'0.5947"
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("D39").Value = "'0.01666"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("D40").Value = "'5.877"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "'0.8488"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Value = "1.047.44"
$ws.Range("E42").Value = "  +2.68%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'101.81"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "1.882.80"
$ws.Range("E45").Value = "  +5.16%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +13.17%  "
$ws.Range("D47").Value = "'59.17"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'8.254"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").Value = "'0.9997"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'0.05315"
$ws.Range("E51").Value = "  -0.74%  "
